# Fixed inspection plan logic and column mapping
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("C8").Value = "MLK_PMT_10101_-_V-001"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "Shell"
$ws.Range("G8").Value = "WATER / STEAM"
$ws.Range("H8").Value = "Carbon Steel"
$ws.Range("I8").Value = "SA-516"
$ws.Range("J8").Value = "'70"
$ws.Range("L8").Value = "80 °C"
$ws.Range("M8").Value = "14 BAR"
$ws.Range("N8").Value = "70 °C"
$ws.Range("O8").Value = "12 BAR"

# --- Row 9 ---
$ws.Range("E9").Value = "Top Head"
$ws.Range("G9").Value = "WATER / STEAM"
$ws.Range("H9").Value = "Carbon Steel"
$ws.Range("I9").Value = "SA-516"
$ws.Range("J9").Value = "'70"
$ws.Range("L9").Value = "80 °C"
$ws.Range("M9").Value = "14 BAR"
$ws.Range("N9").Value = "70 °C"
$ws.Range("O9").Value = "12 BAR"

# --- Row 10 ---
$ws.Range("E10").Value = "Bottom Head"
$ws.Range("G10").Value = "WATER / STEAM"
$ws.Range("H10").Value = "Carbon Steel"
$ws.Range("I10").Value = "SA-516"
$ws.Range("J10").Value = "'70"
$ws.Range("L10").Value = "80 °C"
$ws.Range("M10").Value = "14 BAR"
$ws.Range("N10").Value = "70 °C"
$ws.Range("O10").Value = "12 BAR"

# --- Row 11 ---
$ws.Range("E11").Value = "Flange"
$ws.Range("G11").Value = "WATER / STEAM"
$ws.Range("H11").Value = "Carbon Steel"
$ws.Range("I11").Value = "SA-105"
$ws.Range("J11").Value = ""
$ws.Range("L11").Value = "80 °C"
$ws.Range("M11").Value = "14 BAR"
$ws.Range("N11").Value = "70 °C"
$ws.Range("O11").Value = "12 BAR"

# --- Row 12 ---
$ws.Range("E12").Value = "Lifting Lug"
$ws.Range("G12").Value = "WATER / STEAM"
$ws.Range("H12").Value = "Carbon Steel"
$ws.Range("I12").Value = "SA-516"
$ws.Range("J12").Value = "'70"
$ws.Range("L12").Value = "80 °C"
$ws.Range("M12").Value = "14 BAR"
$ws.Range("N12").Value = "70 °C"
$ws.Range("O12").Value = "12 BAR"

# --- Row 13 ---
$ws.Range("E13").Value = "Name Plate"
$ws.Range("G13").Value = "WATER / STEAM"
$ws.Range("I13").Value = "SS"
$ws.Range("J13").Value = "'304"
$ws.Range("L13").Value = "80 °C"
$ws.Range("M13").Value = "14 BAR"
$ws.Range("N13").Value = "70 °C"
$ws.Range("O13").Value = "12 BAR"

# --- Row 14 ---
$ws.Range("E14").Value = "Bolt & Nut"
$ws.Range("G14").Value = "WATER / STEAM"
$ws.Range("H14").Value = "Stainless Steel Bolting"
$ws.Range("I14").Value = "SA-193 GR. B7 / SA-194 GR. 2H"
$ws.Range("J14").Value = ""
$ws.Range("L14").Value = "80 °C"
$ws.Range("M14").Value = "14 BAR"
$ws.Range("N14").Value = "70 °C"
$ws.Range("O14").Value = "12 BAR"

# --- Row 15 ---
$ws.Range("E15").Value = "Leg Support"
$ws.Range("G15").Value = "WATER / STEAM"
$ws.Range("H15").Value = "Carbon Steel"
$ws.Range("I15").Value = "SA-516"
$ws.Range("J15").Value = "'70"
$ws.Range("L15").Value = "80 °C"
$ws.Range("M15").Value = "14 BAR"
$ws.Range("N15").Value = "70 °C"
$ws.Range("O15").Value = "12 BAR"

# --- Row 16 ---
$ws.Range("E16").Value = "Pipe Clip"
$ws.Range("G16").Value = "WATER / STEAM"
$ws.Range("H16").Value = "Carbon Steel"
$ws.Range("I16").Value = "SA-516"
$ws.Range("J16").Value = "'70"
$ws.Range("L16").Value = "80 °C"
$ws.Range("M16").Value = "14 BAR"
$ws.Range("N16").Value = "70 °C"
$ws.Range("O16").Value = "12 BAR"

# --- Remove now-obsolete rows 17-21 (also shrinks merged ranges A8:A21,
# B8:B21, C8:C21, D8:D21 down to ...:16 and fixes the sheet dimension) ---
$ws.Rows("17:21").Delete()
